$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.076.18"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.68%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.731.28"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.31%  "

$ws.Range("E4").Value = "  -0.19%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.03%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.00%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4878"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +6.48%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3513"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.67%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "43.58"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.99%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07278"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.87%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.051"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.39%  "

$ws.Range("E12").Value = "  -0.13%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.03"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.25%  "

$ws.Range("E14").Value = "  -1.06%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.720.50"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.94%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.900"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.27%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "87.29"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.75%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001040"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.10%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06413"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.13%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.00%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.62"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.03%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.702"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.80%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.126.10"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.68%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.91"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.67%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.080"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.67%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.22"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.86%  "

$ws.Range("E27").Value = "  +0.14%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.929.16"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.43%  "

$ws.Range("E29").Value = "  -1.95%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "121.45"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.00%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.053"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.74%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09399"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.70%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.659"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.09%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.395"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.24%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.06004"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.16%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02195"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.87%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.440"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.89%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "11.01"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.50%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.786"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.09%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2002"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.69%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6020"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.01%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.0000"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.03%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.103"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.26%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.520"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.20%  "

$ws.Range("E45").Value = "  -2.15%  "

$ws.Range("E46").Value = "  -3.65%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5664"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.56%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "119.19"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.93%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.856"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.33%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.104"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.28%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06655"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.11%  "
